# "Análise das vogais /o/ e /u/"
# Fill in the recording data for row 7 (Vogal /o/) and row 8 (Vogal /u/)
# that was previously missing (row 8 triggered #DIV/0! errors).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Individual")

# --- Row 7: Vogal /o/ ---
$ws.Range("B7").Value = 414.8
$ws.Range("C7").Value = 1038
$ws.Range("D7").Value = 2273
$ws.Range("E7").Value = 3486.5
$ws.Range("F7").Value = 404.2
$ws.Range("G7").Value = 1012
$ws.Range("H7").Value = 2329
$ws.Range("I7").Value = 3906
$ws.Range("J7").Value = 423.2
$ws.Range("K7").Value = 846.8
$ws.Range("L7").Value = 2326
$ws.Range("M7").Value = 3067

# --- Row 8: Vogal /u/ ---
$ws.Range("B8").Value = 118.2
$ws.Range("C8").Value = 709.4
$ws.Range("D8").Value = 2943
$ws.Range("E8").Value = 3906
$ws.Range("F8").Value = 119.8
$ws.Range("G8").Value = 831.6
$ws.Range("H8").Value = 2489
$ws.Range("I8").Value = 3962
$ws.Range("J8").Value = 124.4
$ws.Range("K8").Value = 751.9
$ws.Range("L8").Value = 1953
$ws.Range("M8").Value = 3138

# O8 previously had a stray underline font (style for the #DIV/0! cell);
# now that it holds a normal average value, drop the underline.
$ws.Range("O8").Font.Underline = [Microsoft.Office.Interop.Excel.XlUnderlineStyle]::xlUnderlineStyleNone

# Reposition/resize the chart (it was nudged/resized slightly).
$co = $ws.ChartObjects(1)
$co.Left = 371.8927103838583
$co.Top = 141.37393700787402
$co.Width = 432.3383267716535
$co.Height = 216.0

# Update the active selection on the sheet.
$ws.Range("G9").Select()
